# Append: 2025-11-18 18:34 JST
# Update the "取得日時" (acquisition timestamp) column on the ランサーズ sheet
# for every existing data row (rows 2-14) to the refreshed timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-11-18 18:34:57"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 14 }

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}
